$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4 and 5 down to 5 and 6
$ws.Rows("4:4").Insert()

# Fill in the new row 4 data
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 61
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = "train_dim1_2"

# Fix up A column sequential numbering for rows 5 and 6 (A5=4, A6=5)
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Update I column (praclen) to 5 for all data rows
$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 5

# Update sheet view selection to match diff (I6 -> I7)
$ws.Range("I7").Select()
